$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# Row 7 : aluminum tube 1.25"OD
# --------------------------------------------------------------------------
$ws.Range('B7').Value = 'aluminum tube 1.25"OD'
$ws.Range('H7').Value = '1.12"ID'
$ws.Range('C7').Value = 'N'
$ws.Range('D7').Value = 'https://www.mcmaster.com/#9056k76/=19bq239'
$ws.Range('E7').Value = 12.1
$ws.Range('F7').Value = 1
$ws.Range('G7').Value = 'feet'
$ws.Range('L7').Formula = '=E7*F7*I7*(1+J7)'
$ws.Range('M7').Formula = '=M6+L7'

# --------------------------------------------------------------------------
# Row 8 : aluminum rod 1.25"OD
# --------------------------------------------------------------------------
$ws.Range('D8').Value = 'https://www.mcmaster.com/#8974k16/=19bq4f4'
$ws.Range('B8').Value = 'aluminum rod 1.25"OD'
$ws.Range('C8').Value = 'N'
$ws.Range('E8').Value = 11.04
$ws.Range('F8').Value = 1
$ws.Range('G8').Value = 'feet'
$ws.Range('I8').Value = 1
$ws.Range('J8').Value = 0.08
$ws.Range('K8').Value = 0
$ws.Range('L8').Formula = '=E8*F8*I8*(1+J8)'
$ws.Range('M8').Formula = '=M7+L8'

# --------------------------------------------------------------------------
# Row 12 : copper disk 3/4" (filled in out of order, before row 9/10/13)
# --------------------------------------------------------------------------
$ws.Range('H12').Value = '20 guage'
$ws.Range('D12').Value = 'https://www.etsy.com/listing/62082177/copper-discs-20-gauge-stamping-blanks'

# --------------------------------------------------------------------------
# Row 11 : Electrical Connects (section header)
# --------------------------------------------------------------------------
$ws.Range('A11').Value = 'Electrical Connects'

# --------------------------------------------------------------------------
# back to row 12
# --------------------------------------------------------------------------
$ws.Range('B12').Value = 'copper disk 3/4"'
$ws.Range('C12').Value = 'N'
$ws.Range('E12').Formula = '=4/10'
$ws.Range('F12').Value = 70
$ws.Range('G12').Value = 'units'
$ws.Range('I12').Value = 1
$ws.Range('J12').Value = 0.08
$ws.Range('K12').Value = 4.6
$ws.Range('L12').Formula = '=E12*F12*I12*(1+J12)'
$ws.Range('M12').Formula = '=M11+L12'

# --------------------------------------------------------------------------
# Row 13 : copper braid
# --------------------------------------------------------------------------
$ws.Range('D13').Value = 'https://www.amazon.com/dp/B003HGHQVU/ref=biss_dp_t_asn'
$ws.Range('B13').Value = 'copper braid'
$ws.Range('G13').Value = 'roll'
$ws.Range('H13').Value = '25''x.25"'
$ws.Range('C13').Value = 'N'
$ws.Range('E13').Value = 20.01
$ws.Range('F13').Value = 1
$ws.Range('I13').Value = 1
$ws.Range('J13').Value = 0.08
$ws.Range('K13').Value = 0
$ws.Range('L13').Formula = '=E13*F13*I13*(1+J13)'
$ws.Range('M13').Formula = '=M12+L13'

# --------------------------------------------------------------------------
# Row 9 : compression spring
# --------------------------------------------------------------------------
$ws.Range('B9').Value = 'compression spring'
$ws.Range('D9').Value = 'https://www.mcmaster.com/#9657k314/=19bqcva'
$ws.Range('H9').Value = '1"Long .49"ID'
$ws.Range('C9').Value = 'N'
$ws.Range('E9').Formula = '=10.35/12'
$ws.Range('F9').Value = 12
$ws.Range('G9').Value = 'units'
$ws.Range('I9').Value = 1
$ws.Range('J9').Value = 0.08
$ws.Range('K9').Value = 0
$ws.Range('L9').Formula = '=E9*F9*I9*(1+J9)'
$ws.Range('M9').Formula = '=M8+L9'

# --------------------------------------------------------------------------
# Row 14 : Ring Terminal
# --------------------------------------------------------------------------
$ws.Range('B14').Value = 'Ring Terminal'
$ws.Range('H14').Value = '.25"IDx.49"OD'
$ws.Range('D14').Value = 'https://www.mcmaster.com/#7113k444/=19bqfh1'
$ws.Range('C14').Value = 'N'
$ws.Range('E14').Formula = '=12.92/25'
$ws.Range('F14').Value = 25
$ws.Range('G14').Value = 'units'
$ws.Range('I14').Value = 1
$ws.Range('J14').Value = 0.08
$ws.Range('K14').Value = 0
$ws.Range('L14').Formula = '=E14*F14*I14*(1+J14)'
$ws.Range('M14').Formula = '=M13+L14'

# --------------------------------------------------------------------------
# Row 15 : Copper Terminal screw .25"
# --------------------------------------------------------------------------
$ws.Range('B15').Value = 'Copper Terminal screw .25"'
$ws.Range('D15').Value = 'https://www.mcmaster.com/#92949a832/=19bqg8g'
$ws.Range('H15').Value = '1/4"-20, 1/4"length'
$ws.Range('C15').Value = 'N'
$ws.Range('E15').Formula = '=7.03/50'
$ws.Range('F15').Value = 50
$ws.Range('G15').Value = 'units'
$ws.Range('I15').Value = 1
$ws.Range('J15').Value = 0.08
$ws.Range('K15').Value = 0
$ws.Range('L15').Formula = '=E15*F15*I15*(1+J15)'
$ws.Range('M15').Formula = '=M14+L15'

# --------------------------------------------------------------------------
# Row 10 : Black Oxide Screw .375"
# --------------------------------------------------------------------------
$ws.Range('D10').Value = 'https://www.mcmaster.com/#91251a146/=19bqn03'
$ws.Range('B10').Value = 'Black Oxide Screw .375"'
$ws.Range('C10').Value = 'N'
$ws.Range('E10').Formula = '=8.42/100'
$ws.Range('F10').Value = 100
$ws.Range('G10').Value = 'units'
$ws.Range('H10').Value = '6-32 thread'
$ws.Range('I10').Value = 1
$ws.Range('J10').Value = 0.08
$ws.Range('K10').Value = 0
$ws.Range('L10').Formula = '=E10*F10*I10*(1+J10)'
$ws.Range('M10').Formula = '=M9+L10'

# --------------------------------------------------------------------------
# Row 11 numeric cells (I/L/M)
# --------------------------------------------------------------------------
$ws.Range('I11').Value = 1
$ws.Range('L11').Formula = '=E11*F11*I11*(1+J11)'
$ws.Range('M11').Formula = '=M10+L11'

# --------------------------------------------------------------------------
# Hyperlinks - added in the same order the author inserted them
# --------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range('D7'), 'https://www.mcmaster.com/', '9056k76/=19bq239') | Out-Null
$ws.Range('D7').Style = 'Hyperlink'

$ws.Hyperlinks.Add($ws.Range('D8'), 'https://www.mcmaster.com/', '8974k16/=19bq4f4') | Out-Null
$ws.Range('D8').Style = 'Hyperlink'

$ws.Hyperlinks.Add($ws.Range('D12'), 'https://www.etsy.com/listing/62082177/copper-discs-20-gauge-stamping-blanks') | Out-Null
$ws.Range('D12').Style = 'Hyperlink'

$ws.Hyperlinks.Add($ws.Range('D9'), 'https://www.mcmaster.com/', '9657k314/=19bqcva') | Out-Null
$ws.Range('D9').Style = 'Hyperlink'

$ws.Hyperlinks.Add($ws.Range('D14'), 'https://www.mcmaster.com/', '7113k444/=19bqfh1') | Out-Null
$ws.Range('D14').Style = 'Hyperlink'

$ws.Hyperlinks.Add($ws.Range('D15'), 'https://www.mcmaster.com/', '92949a832/=19bqg8g') | Out-Null
$ws.Range('D15').Style = 'Hyperlink'

$ws.Hyperlinks.Add($ws.Range('D10'), 'https://www.mcmaster.com/', '91251a146/=19bqn03') | Out-Null
$ws.Range('D10').Style = 'Hyperlink'

# --------------------------------------------------------------------------
# Final selection (matches the author's last click before saving)
# --------------------------------------------------------------------------
$ws.Range('M16').Select() | Out-Null
